# GPLIM-2980: update the control-reagents test fixture so the third control
# row references a different lot ("SK-2345" instead of "SK-1234").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 is the NO_TEMPLATE_CONTROL row; column C holds the lot barcode.
$ws.Range("C4").Value = "SK-2345"
